$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

$lastRow = $ws.UsedRange.Rows.Count

# Insert a new column before column B, shifting existing B:E to C:F
$ws.Columns("B:B").Insert()

# Header row gets "budget-type"
$ws.Range("B1").Value = "budget-type"

# All data rows (2..lastRow) get "budget"
$ws.Range("B2:B" + $lastRow).Value = "budget"
